$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.198.43"
$ws.Range("E2").Value = "  +0.38%  "
$ws.Range("D3").Value = "2.624.70"
$ws.Range("E3").Value = "  +0.70%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "597.43"
$ws.Range("E5").Value = "  +1.00%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "152.17"
$ws.Range("E6").Value = "  -0.90%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  +2.90%  "
$ws.Range("D9").Value = "2.624.16"
$ws.Range("E9").Value = "  +0.79%  "
$ws.Range("E10").Value = "  +2.59%  "
$ws.Range("E11").Value = "  +0.67%  "
$ws.Range("E12").Value = "  -0.48%  "
$ws.Range("E13").Value = "  -1.26%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.51"
$ws.Range("E14").Value = "  +0.57%  "
$ws.Range("D15").Value = "3.103.56"
$ws.Range("E15").Value = "  +0.60%  "
$ws.Range("E16").Value = "  +0.58%  "
$ws.Range("D17").Value = "67.187.85"
$ws.Range("E17").Value = "  +0.34%  "
$ws.Range("D18").Value = "2.624.82"
$ws.Range("E18").Value = "  +0.70%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.13"
$ws.Range("E19").Value = "  -0.56%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "363.13"
$ws.Range("E20").Value = "  +2.17%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.47"
$ws.Range("E21").Value = "  -3.41%  "
$ws.Range("E22").Value = "  -0.37%  "
$ws.Range("E23").Value = "  +2.56%  "
$ws.Range("E24").Value = "  +0.02%  "
$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "70.78"
$ws.Range("E25").Value = "  +4.91%  "
$ws.Range("B26").Value = "Aptos"
$ws.Range("C26").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.02"
$ws.Range("E26").Value = "  -2.03%  "
$ws.Range("B27").Value = "Binance-PegBSC-USD"
$ws.Range("C27").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.03"
$ws.Range("E27").Value = "  +3.61%  "
$ws.Range("B28").Value = "WrappedeETH"
$ws.Range("C28").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D28").Value = "2.760.08"
$ws.Range("E28").Value = "  +0.59%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0000102"
$ws.Range("E29").Value = "  +0.11%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "571.45"
$ws.Range("E30").Value = "  -6.59%  "
$ws.Range("E31").Value = "  -3.54%  "
$ws.Range("E32").Value = "  -1.72%  "
$ws.Range("E33").Value = "  -0.25%  "
$ws.Range("E34").Value = "  +0.09%  "
$ws.Range("E35").Value = "  -3.76%  "
$ws.Range("E36").Value = "  -1.57%  "
$ws.Range("E37").Value = "  -1.37%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "157.08"
$ws.Range("E38").Value = "  +1.68%  "
$ws.Range("E39").Value = "  -0.43%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.368"
$ws.Range("E40").Value = "  +0.02%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.24"
$ws.Range("E41").Value = "  -3.18%  "
$ws.Range("E42").Value = "  -0.32%  "
$ws.Range("E43").Value = "  -0.26%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "41.19"
$ws.Range("E44").Value = "  -0.18%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.00"
$ws.Range("E45").Value = "  +0.02%  "
$ws.Range("E46").Value = "  -0.62%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "155.60"
$ws.Range("E47").Value = "  +0.46%  "
$ws.Range("E48").Value = "  -2.16%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.73"
$ws.Range("E50").Value = "  -0.36%  "
$ws.Range("E51").Value = "  -1.51%  "

# Restore default style for cells that were temporarily forced to Text format,
# so they match the original (un-styled) appearance while keeping string content.
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D49").Style = "Normal"
